$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.303.52"
$ws.Range("E2").Value = "  +2.75%  "
$ws.Range("D3").Value = "3.490.24"
$ws.Range("E3").Value = "  +2.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.21%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.480"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("E10").Value = "  +3.12%  "
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("D12").Value = "4.085.87"
$ws.Range("E12").Value = "  +2.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.06%  "
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").Value = "3.488.92"
$ws.Range("E15").Value = "  +2.96%  "
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").Value = "63.304.20"
$ws.Range("E17").Value = "  +2.85%  "
$ws.Range("E18").Value = "  +2.32%  "
$ws.Range("E19").Value = "  +5.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "391.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.567"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "75.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000119"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.15%  "
$ws.Range("D26").Value = "3.629.08"
$ws.Range("E26").Value = "  +2.93%  "
$ws.Range("E27").Value = "  -4.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.64%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +3.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.75%  "
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("E35").Value = "  +6.32%  "
$ws.Range("E36").Value = "  +2.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "32.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +23.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "171.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("E39").Value = "  +6.67%  "
$ws.Range("D40").Value = "3.525.33"
$ws.Range("E40").Value = "  +2.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0772"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  +3.98%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.99%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("E46").Value = "  +7.36%  "
$ws.Range("D47").Value = "2.626.25"
$ws.Range("E47").Value = "  +6.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.12%  "
$ws.Range("E49").Value = "  +12.61%  "
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("E51").Value = "  +3.43%  "
